# Changing the Risk-free rate to zero: update the computed optimal
# portfolio weights (columns C and D) for rows 2-8 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.562341192568151
$ws.Range("D2").Value = 0.562341192568151

$ws.Range("C3").Value = 0.3584989408622693
$ws.Range("D3").Value = 0.3584989408622693

$tinyVal = [double]"4.336808689942018e-19"
$ws.Range("C4").Value = $tinyVal
$ws.Range("D4").Value = $tinyVal

$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("C8").Value = 0.07915986656957973
$ws.Range("D8").Value = 0.07915986656957973
